$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group_a")

# Proxy moved to a new machine: update PROXY:PORT, PROXY_USER, PROXY_PASS
# for the John Snow row (row 2) to the new values.
$ws.Range("E2").Value = "81.28.96.148:4000"
$ws.Range("F2").Value = "i0BdGW79w6Oo"
$ws.Range("G2").Value = "5Ao37R1ry6bc"

$ws.Range("G2").Select()
